$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 6).Value = 3   # F2: 2.96 -> 3
$ws.Cells.Item(2, 8).Value = 2.42   # H2: 2.44 -> 2.42
$ws.Cells.Item(2, 9).Value = 2.7   # I2: 2.72 -> 2.7
$ws.Cells.Item(2, 10).Value = 3.15   # J2: 3.1 -> 3.15
$ws.Cells.Item(2, 12).Value = 1.45   # L2: 1.44 -> 1.45
$ws.Cells.Item(2, 17).Value = 2.1   # Q2: 2.08 -> 2.1
$ws.Cells.Item(2, 22).Value = 1.59   # V2: 1.58 -> 1.59
$ws.Cells.Item(3, 8).Value = 2.7   # H3: 2.66 -> 2.7
$ws.Cells.Item(3, 11).Value = 4.1   # K3: 4.2 -> 4.1
$ws.Cells.Item(3, 20).Value = 1.71   # T3: 1.7 -> 1.71
$ws.Cells.Item(3, 25).Value = 14.5   # Y3: 15 -> 14.5
$ws.Cells.Item(4, 8).Value = 2.82   # H4: 2.8 -> 2.82
$ws.Cells.Item(4, 9).Value = 3.25   # I4: 3.2 -> 3.25
$ws.Cells.Item(4, 10).Value = 3.3   # J4: 3.15 -> 3.3
$ws.Cells.Item(4, 11).Value = 3.85   # K4: 3.9 -> 3.85
$ws.Cells.Item(4, 12).Value = 1.39   # L4: 1.33 -> 1.39
$ws.Cells.Item(4, 14).Value = 3.6   # N4: 3.7 -> 3.6
$ws.Cells.Item(4, 16).Value = 1.9   # P4: 1.93 -> 1.9
$ws.Cells.Item(4, 18).Value = 1.34   # R4: 1.36 -> 1.34
$ws.Cells.Item(4, 24).Value = 17.5   # X4: 1000 -> 17.5
$ws.Cells.Item(4, 29).Value = 8.800000000000001   # AC4: 1000 -> 8.800000000000001
$ws.Cells.Item(5, 6).Value = 8.6   # F5: 8.4 -> 8.6
$ws.Cells.Item(5, 7).Value = 8.800000000000001   # G5: 8.6 -> 8.800000000000001
$ws.Cells.Item(5, 9).Value = 1.44   # I5: 1.45 -> 1.44
$ws.Cells.Item(5, 11).Value = 5.5   # K5: 5.4 -> 5.5
$ws.Cells.Item(5, 12).Value = 1.32   # L5: 1.33 -> 1.32
$ws.Cells.Item(5, 14).Value = 5.3   # N5: 5.2 -> 5.3
$ws.Cells.Item(5, 17).Value = 1.65   # Q5: 1.66 -> 1.65
$ws.Cells.Item(5, 18).Value = 1.57   # R5: 1.56 -> 1.57
$ws.Cells.Item(5, 19).Value = 2.68   # S5: 2.72 -> 2.68
$ws.Cells.Item(5, 22).Value = 3.25   # V5: 3.2 -> 3.25
$ws.Cells.Item(5, 23).Value = 1.12   # W5: 1.13 -> 1.12
$ws.Cells.Item(5, 24).Value = 24   # X5: 23 -> 24
$ws.Cells.Item(5, 25).Value = 9.800000000000001   # Y5: 9.6 -> 9.800000000000001
$ws.Cells.Item(5, 34).Value = 23   # AH5: 25 -> 23
$ws.Cells.Item(5, 35).Value = 30   # AI5: 32 -> 30
$ws.Cells.Item(5, 38).Value = 95   # AL5: 100 -> 95
$ws.Cells.Item(5, 39).Value = 120   # AM5: 130 -> 120
$ws.Cells.Item(5, 41).Value = 5.6   # AO5: 5.8 -> 5.6
$ws.Cells.Item(6, 6).Value = 4.6   # F6: 4.7 -> 4.6
$ws.Cells.Item(6, 11).Value = 4.5   # K6: 4.4 -> 4.5
$ws.Cells.Item(6, 14).Value = 5.8   # N6: 5.7 -> 5.8
$ws.Cells.Item(6, 32).Value = 980   # AF6: 46 -> 980
$ws.Cells.Item(7, 18).Value = 1.49   # R7: 1.5 -> 1.49
$ws.Cells.Item(7, 19).Value = 2.96   # S7: 2.94 -> 2.96
$ws.Cells.Item(8, 14).Value = 1.24   # N8: 1.1 -> 1.24
$ws.Cells.Item(8, 16).Value = 1.24   # P8: 1.09 -> 1.24
$ws.Cells.Item(8, 18).Value = 1.18   # R8: 1.08 -> 1.18
$ws.Cells.Item(8, 19).Value = 1.3   # S8: 1.27 -> 1.3
$ws.Cells.Item(9, 6).Value = 2.36   # F9: 2.34 -> 2.36
$ws.Cells.Item(9, 22).Value = 1.4   # V9: 1.39 -> 1.4
$ws.Cells.Item(9, 28).Value = 9.800000000000001   # AB9: 10 -> 9.800000000000001
$ws.Cells.Item(10, 17).Value = 1.75   # Q10: 1.74 -> 1.75
$ws.Cells.Item(10, 41).Value = 5.1   # AO10: 5.2 -> 5.1
$ws.Cells.Item(11, 6).Value = 2.68   # F11: 2.72 -> 2.68
$ws.Cells.Item(11, 7).Value = 2.7   # G11: 2.74 -> 2.7
$ws.Cells.Item(11, 8).Value = 2.9   # H11: 2.88 -> 2.9
$ws.Cells.Item(11, 9).Value = 2.94   # I11: 2.92 -> 2.94
$ws.Cells.Item(11, 16).Value = 2   # P11: 1.99 -> 2
$ws.Cells.Item(11, 17).Value = 1.98   # Q11: 1.99 -> 1.98
$ws.Cells.Item(11, 22).Value = 1.51   # V11: 1.52 -> 1.51
$ws.Cells.Item(11, 23).Value = 1.58   # W11: 1.57 -> 1.58
$ws.Cells.Item(11, 24).Value = 13.5   # X11: 14 -> 13.5
$ws.Cells.Item(11, 26).Value = 19   # Z11: 18.5 -> 19
$ws.Cells.Item(11, 29).Value = 7.4   # AC11: 7.6 -> 7.4
$ws.Cells.Item(11, 31).Value = 32   # AE11: 30 -> 32
$ws.Cells.Item(11, 37).Value = 28   # AK11: 29 -> 28
$ws.Cells.Item(12, 6).Value = 4.2   # F12: 4.3 -> 4.2
$ws.Cells.Item(12, 8).Value = 1.93   # H12: 1.91 -> 1.93
$ws.Cells.Item(12, 9).Value = 1.94   # I12: 1.92 -> 1.94
$ws.Cells.Item(12, 12).Value = 1.34   # L12: 1.33 -> 1.34
$ws.Cells.Item(12, 17).Value = 1.71   # Q12: 1.69 -> 1.71
$ws.Cells.Item(12, 22).Value = 2.06   # V12: 2.08 -> 2.06
$ws.Cells.Item(12, 23).Value = 1.3   # W12: 1.29 -> 1.3
$ws.Cells.Item(12, 24).Value = 19   # X12: 19.5 -> 19
$ws.Cells.Item(12, 31).Value = 17.5   # AE12: 17 -> 17.5
$ws.Cells.Item(12, 34).Value = 16   # AH12: 16.5 -> 16
$ws.Cells.Item(12, 36).Value = 85   # AJ12: 90 -> 85
$ws.Cells.Item(12, 37).Value = 44   # AK12: 46 -> 44
$ws.Cells.Item(12, 38).Value = 46   # AL12: 48 -> 46
$ws.Cells.Item(13, 6).Value = 4.7   # F13: 4.8 -> 4.7
$ws.Cells.Item(13, 7).Value = 4.8   # G13: 4.9 -> 4.8
$ws.Cells.Item(13, 17).Value = 1.6   # Q13: 1.59 -> 1.6
$ws.Cells.Item(13, 22).Value = 2.22   # V13: 2.24 -> 2.22
$ws.Cells.Item(13, 24).Value = 22   # X13: 23 -> 22
$ws.Cells.Item(13, 30).Value = 10   # AD13: 9.800000000000001 -> 10
$ws.Cells.Item(13, 36).Value = 95   # AJ13: 100 -> 95
$ws.Cells.Item(13, 38).Value = 46   # AL13: 48 -> 46
$ws.Cells.Item(14, 6).Value = 2.04   # F14: 2.02 -> 2.04
$ws.Cells.Item(14, 8).Value = 3.75   # H14: 4.1 -> 3.75
$ws.Cells.Item(14, 12).Value = 1.43   # L14: 1.01 -> 1.43
$ws.Cells.Item(14, 13).Value = 1.1   # M14: 1.09 -> 1.1
$ws.Cells.Item(14, 16).Value = 1.65   # P14: 1.7 -> 1.65
$ws.Cells.Item(14, 17).Value = 2.28   # Q14: 2.26 -> 2.28
$ws.Cells.Item(14, 18).Value = 1.25   # R14: 1.26 -> 1.25
$ws.Cells.Item(14, 20).Value = 1.98   # T14: 1.96 -> 1.98
$ws.Cells.Item(14, 22).Value = 1.29   # V14: 1.28 -> 1.29
